# Append the latest Adafruit IO reading as a new row at the bottom of the feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 46

$ws.Range("A$newRow").Value = "2024-09-25T18:06:40Z"
$ws.Range("B$newRow").Value = "temperature"
# "Value" column holds numeric-looking readings but the sheet stores every
# column as text, so force a text value (leading apostrophe) instead of
# letting Excel auto-coerce "25" into a number.
$ws.Range("C$newRow").Value = "'25"
$ws.Range("D$newRow").Value = "N/A"
$ws.Range("E$newRow").Value = "N/A"
$ws.Range("F$newRow").Value = "N/A"
